# Adds a new "alias" column to the "Experiências" and "Formações" sheets,
# giving every entry a short machine-friendly slug used elsewhere in the
# curriculum-building pipeline, and updates the active sheet/selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Experiências" sheet: insert a new column C ("alias") and shift
#    cargo_pt/empresa_pt/... etc. one column to the right (C:K -> D:L).
# ---------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("Experiências")
$wsExp.Columns("C:C").Insert()

$wsExp.Range("C1").Value = "alias"

$wsExp.Range("C2").Value = "nov"
$wsExp.Range("C3").Value = "luna"
$wsExp.Range("C4").Value = "aqn"
$wsExp.Range("C5").Value = "eco"
$wsExp.Range("C6").Value = "yop"
$wsExp.Range("C7").Value = "pipa"

# ---------------------------------------------------------------
# 2) "Formações" sheet: insert a new column B ("alias") and shift
#    curso_pt/instituicao_pt/... etc. one column to the right (B:K -> C:L).
# ---------------------------------------------------------------
$wsForm = $wb.Worksheets.Item("Formações")
$wsForm.Columns("B:B").Insert()
$wsForm.Columns("B:B").ColumnWidth = 14.3

$wsForm.Range("B1").Value = "alias"

$wsForm.Range("B2").Value = "grad"
$wsForm.Range("B3").Value = "mic"
$wsForm.Range("B4").Value = "fer"
$wsForm.Range("B5").Value = "csap"
$wsForm.Range("B6").Value = "agil"

# ---------------------------------------------------------------
# 3) Selection / active-sheet bookkeeping, matching the saved workbook
#    state: "Formações" keeps a remembered selection on B7, while
#    "Experiências" becomes the active tab with C1 selected.
# ---------------------------------------------------------------
$wsForm.Activate()
$wsForm.Range("B7").Select()

$wsExp.Activate()
$wsExp.Range("C1").Select()
